$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.747665643692017
$ws.Range("B1").Value = 2.242830991744995
$ws.Range("C1").Value = 2.402452230453491
$ws.Range("D1").Value = 7.359123229980469
$ws.Range("E1").Value = 0.7696225643157959
